$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert 20 new blank rows (bottom-up) to expand each Movies/Music pair into a Movies/Music/Sport triplet
$insertRows = @(622,620,618,616,614,612,610,608,606,604,602,600,598,596,594,592,590,588,586,584)
foreach ($r in $insertRows) {
    $ws.Rows.Item($r).Insert()
}

# Step 2: fill in the newly inserted rows with the "Sport" domain content
$ws.Range("A584").Value = 'Sport & Enterteinment'
$ws.Range("B584").Value = 8
$ws.Range("C584").Value = 'Sport'
$ws.Range("D584").Value = 'adversarial learning'
$ws.Range("E584").Value = 'As a coach, I want to utilize adversarial learning models to simulate opponents'' strategies and tactics based on historical data, in order to better prepare my team for upcoming matches.'
$ws.Range("F584").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A587").Value = 'Sport & Enterteinment'
$ws.Range("B587").Value = 8
$ws.Range("C587").Value = 'Sport'
$ws.Range("D587").Value = 'cnn'
$ws.Range("E587").Value = 'As a sports video analyst, I want to use CNNs to automatically track and analyze player movements during games, providing coaches with detailed insights into player positioning and performance metrics.'
$ws.Range("F587").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A590").Value = 'Sport & Enterteinment'
$ws.Range("B590").Value = 8
$ws.Range("C590").Value = 'Sport'
$ws.Range("D590").Value = 'conversational agent'
$ws.Range("E590").Value = 'As a sports team manager, I want to develop a conversational agent powered by machine learning to provide real-time updates on player injuries, performance statistics, and training schedules, enhancing communication and coordination within the team.'
$ws.Range("F590").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A593").Value = 'Sport & Enterteinment'
$ws.Range("B593").Value = 8
$ws.Range("C593").Value = 'Sport'
$ws.Range("D593").Value = 'decision tree'
$ws.Range("E593").Value = 'As a sports physiologist, I want to employ a decision tree algorithm to classify and diagnose common sports injuries based on symptoms, patient history, and diagnostic tests, aiding in timely and accurate treatment plans.'
$ws.Range("F593").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A596").Value = 'Sport & Enterteinment'
$ws.Range("B596").Value = 8
$ws.Range("C596").Value = 'Sport'
$ws.Range("D596").Value = 'document classification'
$ws.Range("E596").Value = 'As a sports journalist, I want to develop a document classification model to automatically categorize news articles into relevant sports categories such as football, basketball, and tennis, facilitating quicker content retrieval and publication.'
$ws.Range("F596").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A599").Value = 'Sport & Enterteinment'
$ws.Range("B599").Value = 8
$ws.Range("C599").Value = 'Sport'
$ws.Range("D599").Value = 'entity extraction'
$ws.Range("E599").Value = 'As a sports betting analyst, I want to develop an entity extraction model to extract key statistics such as player performance metrics, injury updates, and historical match results from sports betting websites and databases, facilitating data-driven betting strategies.'
$ws.Range("F599").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A602").Value = 'Sport & Enterteinment'
$ws.Range("B602").Value = 8
$ws.Range("C602").Value = 'Sport'
$ws.Range("D602").Value = 'feature selection'
$ws.Range("E602").Value = 'As a sports performance analyst, I want to employ feature selection techniques to identify the most relevant player performance metrics (such as goals scored, assists, and accuracy) that correlate with team success, aiding in player evaluation and strategy formulation.'
$ws.Range("F602").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A605").Value = 'Sport & Enterteinment'
$ws.Range("B605").Value = 8
$ws.Range("C605").Value = 'Sport'
$ws.Range("D605").Value = 'imbalanced dataset'
$ws.Range("E605").Value = 'As a sports talent scout, I want to build a recruitment model that handles imbalanced data by accurately identifying and prioritizing promising young athletes from underrepresented regions or sports disciplines, ensuring comprehensive talent evaluation'
$ws.Range("F605").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A608").Value = 'Sport & Enterteinment'
$ws.Range("B608").Value = 8
$ws.Range("C608").Value = 'Sport'
$ws.Range("D608").Value = 'keyword extraction'
$ws.Range("E608").Value = 'As a sports content curator, I want to implement keyword extraction techniques to analyze player interviews and press conferences, automatically extracting key themes and quotes for creating engaging multimedia content for fans.'
$ws.Range("F608").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A611").Value = 'Sport & Enterteinment'
$ws.Range("B611").Value = 8
$ws.Range("C611").Value = 'Sport'
$ws.Range("D611").Value = 'k-nearest neighbor'
$ws.Range("E611").Value = 'As a sports talent scout, I want to use k-NN clustering to group young athletes based on their physical attributes, skill levels, and potential for development, facilitating more targeted scouting and recruitment efforts.'
$ws.Range("F611").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A614").Value = 'Sport & Enterteinment'
$ws.Range("B614").Value = 8
$ws.Range("C614").Value = 'Sport'
$ws.Range("D614").Value = 'multi-label classification'
$ws.Range("E614").Value = 'As a sports physiologist, I want to develop a multi-label classification model to predict the physiological responses (such as heart rate, oxygen consumption) of athletes during different phases of a game or training session, based on various environmental and physical factors.'
$ws.Range("F614").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A617").Value = 'Sport & Enterteinment'
$ws.Range("B617").Value = 8
$ws.Range("C617").Value = 'Sport'
$ws.Range("D617").Value = 'neural network'
$ws.Range("E617").Value = 'As a sports video analyst, I want to develop a convolutional neural network (CNN) to automatically analyze player movements and positions from video footage, extracting tactical insights and performance metrics for coaching and strategy refinement.'
$ws.Range("F617").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A620").Value = 'Sport & Enterteinment'
$ws.Range("B620").Value = 8
$ws.Range("C620").Value = 'Sport'
$ws.Range("D620").Value = 'random forest'
$ws.Range("E620").Value = 'As a sports betting strategist, I want to build a random forest model to predict betting odds for various outcomes in sports matches, considering factors such as team form, player statistics, and historical match data, to inform strategic betting decisions.'
$ws.Range("F620").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A623").Value = 'Sport & Enterteinment'
$ws.Range("B623").Value = 8
$ws.Range("C623").Value = 'Sport'
$ws.Range("D623").Value = 'semantic similarity'
$ws.Range("E623").Value = 'As a sports content curator, I want to develop a semantic similarity model to recommend relevant articles, videos, and social media posts to fans based on their interests and engagement history, enhancing personalized content delivery.'
$ws.Range("F623").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A626").Value = 'Sport & Enterteinment'
$ws.Range("B626").Value = 8
$ws.Range("C626").Value = 'Sport'
$ws.Range("D626").Value = 'sentiment analysis'
$ws.Range("E626").Value = 'As a sports broadcaster, I want to develop a sentiment analysis model to analyze viewer reactions and sentiment towards live sports broadcasts, helping us understand audience engagement and preferences.'
$ws.Range("F626").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A629").Value = 'Sport & Enterteinment'
$ws.Range("B629").Value = 8
$ws.Range("C629").Value = 'Sport'
$ws.Range("D629").Value = 'speech to text'
$ws.Range("E629").Value = 'As a sports journalist, I want to leverage speech-to-text technology to transcribe interviews with athletes and coaches, enabling faster content creation and accurate reporting of quotes and insights.'
$ws.Range("F629").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A632").Value = 'Sport & Enterteinment'
$ws.Range("B632").Value = 8
$ws.Range("C632").Value = 'Sport'
$ws.Range("D632").Value = 'text categorization'
$ws.Range("E632").Value = 'As a sports news aggregator, I want to develop a text categorization model to classify news articles into different sports categories such as football, basketball, tennis, etc., ensuring relevant and organized content delivery to users.'
$ws.Range("F632").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A635").Value = 'Sport & Enterteinment'
$ws.Range("B635").Value = 8
$ws.Range("C635").Value = 'Sport'
$ws.Range("D635").Value = 'unsupervised clustering'
$ws.Range("E635").Value = 'As a sports scouting coordinator, I want to use unsupervised clustering techniques to cluster and analyze player statistics and attributes to identify emerging talents and potential recruits, streamlining talent scouting and recruitment processes.'
$ws.Range("F635").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A638").Value = 'Sport & Enterteinment'
$ws.Range("B638").Value = 8
$ws.Range("C638").Value = 'Sport'
$ws.Range("D638").Value = 'voice recognition'
$ws.Range("E638").Value = 'As a sports broadcaster, I want to deploy voice recognition systems to automatically generate transcripts of sports podcasts and radio shows, enabling searchable archives and content summaries for listeners.'
$ws.Range("F638").Value = 'ReAdjusted_CoTPrompt'

$ws.Range("A641").Value = 'Sport & Enterteinment'
$ws.Range("B641").Value = 8
$ws.Range("C641").Value = 'Sport'
$ws.Range("D641").Value = 'word embedding'
$ws.Range("E641").Value = 'As a sports content curator, I want to apply word embedding models to recommend related articles and videos to sports fans based on semantic similarities in content, enhancing user engagement and content discovery.'
$ws.Range("F641").Value = 'ReAdjusted_CoTPrompt'
